{"js": "// The document contains a stray, empty bulleted list paragraph that was\n// left behind right after the \"create threads out of questions\" bullet\n// (and right before the document's trailing empty paragraph). Remove\n// that empty paragraph entirely.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst anchorText = \"The user shall be able to create threads out of questions \";\n\nlet emptyParaIndex = -1;\nfor (let i = 0; i < paragraphs.items.length - 1; i++) {\n  if (paragraphs.items[i].text === anchorText && paragraphs.items[i + 1].text === \"\") {\n    emptyParaIndex = i + 1;\n    break;\n  }\n}\n\nif (emptyParaIndex === -1) {\n  throw new Error(\"Could not find the stray empty paragraph to delete.\");\n}\n\nparagraphs.items[emptyParaIndex].delete();\nawait context.sync();\n", "ps1": "# The document contains a stray, empty bulleted list paragraph that was\n# left behind right after the \"create threads out of questions\" bullet\n# (and right before the document's trailing empty paragraph). Remove\n# that empty paragraph entirely.\n\n$d = $word.ActiveDocument\n\n$anchorText = \"The user shall be able to create threads out of questions \"\n\n$searchRange = $d.Content\n$find = $searchRange.Find\n$find.ClearFormatting()\n$find.Text = $anchorText\n$found = $find.Execute()\n\nif (-not $found) {\n    throw \"Could not find the anchor paragraph text.\"\n}\n\n$anchorPara = $searchRange.Paragraphs(1)\n$emptyPara = $anchorPara.Next()\n\n$emptyText = $emptyPara.Range.Text.TrimEnd([char]13, [char]7)\nif ($emptyText -ne \"\") {\n    throw \"Paragraph following the anchor was not empty; refusing to delete.\"\n}\n\n$emptyPara.Range.Delete()\n"}
